$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Behavior / Knowledge..." sentence.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Behavior / Knowledge: relation tuples rows data*") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph 'Behavior / Knowledge...'"
}

# The document already has an (empty) paragraph directly after the anchor
# paragraph; it is left untouched. All of the new paragraphs below are
# inserted right after that existing blank paragraph, in order. Empty
# strings represent blank paragraphs.
$newParagraphs = @(
    "FCA Augmented Resources.",
    "",
    "Relationship Monad functions: assert.",
    "",
    "Relationship Monad functions: assert.",
    "",
    "ID Relationship Monad:",
    "",
    "Entity Monad ID Relationship:",
    "",
    "anEntity.flatMap(ID::assert(r : Relationship) : s : Relationship (anEntity if equals, previous / next Entity if not equals).",
    ""
)

$existingBlank = $anchor.Next()

# Anchor the insertion point at the *start* of the existing blank paragraph
# (rather than its end) so it stays unambiguously inside that paragraph --
# using the end position sits exactly on the boundary with the following
# paragraph and causes InsertParagraphAfter to insert after the wrong
# paragraph.
$pos = $existingBlank.Range.Start

foreach ($text in $newParagraphs) {
    $r = $d.Range($pos, $pos)
    $r.InsertParagraphAfter()
    if ($text -ne "") {
        $r2 = $d.Range($pos + 1, $pos + 1)
        $r2.InsertAfter($text)
        $pos = $pos + 1 + $text.Length
    } else {
        $pos = $pos + 1
    }
}
